$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# --- Cells that flip between a numeric value and the "0"/"***.*" placeholder text ---
# Copy style+value from a same-column neighbor that already carries the desired
# type/style, then overwrite with the actual target value, so the resulting
# cell keeps the same style index as the rest of the column.

# F15: numeric 1 -> placeholder text "0" (style like G15)
$ws.Range("G15").Copy($ws.Range("F15"))

# C18: placeholder text "0" -> numeric 5 (style like D18)
$ws.Range("D18").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 5

# C20: numeric 2 -> placeholder text "0" (style like C14)
$ws.Range("C14").Copy($ws.Range("C20"))

# D20: numeric 2 -> placeholder text "0" (style like D14)
$ws.Range("D14").Copy($ws.Range("D20"))

# E20: numeric 0 -> placeholder text "***.*" (style like E14)
$ws.Range("E14").Copy($ws.Range("E20"))

# C27: placeholder text "0" -> numeric 1 (style like C28)
$ws.Range("C28").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1

# C29: numeric 1 -> placeholder text "0" (style like D29)
$ws.Range("D29").Copy($ws.Range("C29"))

# C30: numeric 1 -> placeholder text "0" (style like D30)
$ws.Range("D30").Copy($ws.Range("C30"))

# --- Remaining plain numeric value updates ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = -17.525773195876
$ws.Range("L16").Value = -36.507936507936
$ws.Range("M16").Value = 14.285714285714
$ws.Range("N16").Value = -82.142857142857
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -34.285714285714
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 117
$ws.Range("K17").Value = -31.623931623931
$ws.Range("L17").Value = -13.978494623655
$ws.Range("M17").Value = 70.212765957446
$ws.Range("N17").Value = -46.308724832214
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 131
$ws.Range("J18").Value = 164
$ws.Range("K18").Value = -20.121951219512
$ws.Range("L18").Value = -44.255319148936
$ws.Range("M18").Value = 24.761904761904
$ws.Range("N18").Value = -69.320843091334
$ws.Range("C19").Value = 30
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 103
$ws.Range("G19").Value = 110
$ws.Range("H19").Value = -6.363636363636
$ws.Range("I19").Value = 561
$ws.Range("J19").Value = 659
$ws.Range("K19").Value = -14.871016691957
$ws.Range("L19").Value = -14.741641337386
$ws.Range("M19").Value = 1.081081081081
$ws.Range("N19").Value = -55.826771653543
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("L20").Value = -31.25
$ws.Range("N20").Value = -94.387755102040
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 185
$ws.Range("H21").Value = -17.297297297297
$ws.Range("I21").Value = 876
$ws.Range("J21").Value = 1061
$ws.Range("K21").Value = -17.436380772855
$ws.Range("L21").Value = -24.090121317157
$ws.Range("M21").Value = 8.819875776397
$ws.Range("N21").Value = -67.471221685852
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 9.090909090909
$ws.Range("M22").Value = -25
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = -17.777777777777
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 226
$ws.Range("H24").Value = -40.707964601769
$ws.Range("I24").Value = 910
$ws.Range("J24").Value = 1095
$ws.Range("K24").Value = -16.894977168949
$ws.Range("L24").Value = -12.246865959498
$ws.Range("M24").Value = 9.506618531889
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 35
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 97
$ws.Range("G25").Value = 179
$ws.Range("H25").Value = -45.810055865921
$ws.Range("I25").Value = 732
$ws.Range("J25").Value = 825
$ws.Range("K25").Value = -11.272727272727
$ws.Range("L25").Value = -12.857142857142
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -5.128205128205
$ws.Range("I26").Value = 204
$ws.Range("J26").Value = 244
$ws.Range("K26").Value = -16.393442622950
$ws.Range("L26").Value = -2.857142857142
$ws.Range("M26").Value = 53.383458646616
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -50
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 5.555555555555
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 300
$ws.Range("I31").Value = 17
$ws.Range("K31").Value = 183.333333333333
$ws.Range("L31").Value = 112.5
